# Saldo.xlsx update — "Add files via upload"
#
# Net effect (this sheet is kept sorted by Saldo, descending):
#   - Account 008035153 (CLAUDIO)   balance corrected  51.48  ->  28051.48
#   - Account 001368670 (THIAGO)    balance corrected  89.53  ->    189.53
#   - New account 005009922 (ANA)   added with balance 26311.36
# Because the sheet is sorted by the Saldo column, correcting/adding these
# balances means the affected rows move to new positions in the list.
# We implement this as: delete the two old (now-stale) rows, then insert
# three rows (two new + one re-positioned) at their correct sorted slots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, [string]$text) {
    # Force text storage so leading zeros in account numbers survive
    # (plain .Value assignment on an all-digit string is parsed as a number).
    $range.NumberFormat = "@"
    $range.Value = $text
}

function Insert-DataRow([int]$atRow, [string]$conta, [string]$nome, [double]$saldo) {
    $ws.Rows.Item($atRow).Insert()
    Set-TextCell $ws.Cells.Item($atRow, 1) $conta
    Set-TextCell $ws.Cells.Item($atRow, 2) $nome
    $ws.Cells.Item($atRow, 3).Value = $saldo
}

# --- Remove the two stale rows first (processed bottom-up so the row
#     numbers below still refer to the original/current layout). ---

# Old row: 008035153  CLAUDIO   51.48   (row 241)
$ws.Rows.Item(241).Delete()

# Old row: 001368670  THIAGO    89.53   (row 157)
$ws.Rows.Item(157).Delete()

# --- Insert rows at their new sorted positions (still bottom-up). ---

# THIAGO's corrected balance (189.53) belongs right before 005547702 (row 94).
Insert-DataRow 94 "001368670" "THIAGO" 189.53

# CLAUDIO's corrected balance (28051.48) and the new ANA account (26311.36)
# belong right before 005142611 (row 5), in descending-Saldo order.
Insert-DataRow 5 "008035153" "CLAUDIO" 28051.48
Insert-DataRow 6 "005009922" "ANA" 26311.36
